# Duplicate the first paragraph ("Just Testing") as a brand-new line above
# itself, leave a truly blank paragraph after the duplicate, and turn the
# (now third) original paragraph's text into "Thank you very much ".

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs(1)
$firstRange = $firstPara.Range

# Capture the exact OOXML (including the three separate runs: "Just", " ",
# "Testing") of the existing first paragraph so the new paragraph we insert
# is a faithful duplicate, run-for-run and formatting-for-formatting.
$originalXml = $firstRange.WordOpenXML

# Derive a variant of that same OOXML package with the paragraph's runs
# stripped out (keeping its <w:pPr> intact), to use for the new, completely
# blank paragraph.
$emptyXml = $originalXml -replace '(?s)(</w:pPr>).*?(</w:p>)', '$1$2'

# Insert two new empty paragraphs before the original "Just Testing"
# paragraph: the first will receive the duplicated content, the second
# stays blank.
[void]$firstRange.InsertParagraphBefore()
[void]$firstRange.InsertParagraphBefore()

# Fill the brand-new first paragraph with the duplicated "Just Testing"
# content (restoring the original run structure/formatting exactly).
$newPara1 = $d.Paragraphs(1)
[void]$newPara1.Range.InsertXML($originalXml)

# Fill the brand-new second paragraph with the same paragraph formatting but
# no runs at all, leaving it blank.
$newPara2 = $d.Paragraphs(2)
[void]$newPara2.Range.InsertXML($emptyXml)

# The original paragraph is now the 3rd paragraph; replace its text with
# "Thank you very much ".
$targetPara = $d.Paragraphs(3)
[void]$targetPara.Range.Find.Execute("Just Testing", $false, $false, $false, $false, $false, $true, 1, $false, "Thank you very much ", 2)
